$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1763.3784
$ws.Range("I15").Value = 1763.3784
$ws.Range("K15").Value = 5290.135200000001
$ws.Range("M15").Value = -5121.135200000001
$ws.Range("H43").Value = 5107
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H80").Value = 666
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 499
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 1497
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -3493
$ws.Range("H83").Value = 666
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 499
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 4491
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -14475
$ws.Range("H86").Value = 2281.5
$ws.Range("I86").Value = 1633.75
$ws.Range("K86").Value = 1633.75
$ws.Range("M86").Value = -510.75
$ws.Range("H89").Value = 2281.5
$ws.Range("I89").Value = 1633.75
$ws.Range("K89").Value = 8168.75
$ws.Range("M89").Value = -2552.75
$ws.Range("H98").Value = 672.7273
$ws.Range("I98").Value = 853
$ws.Range("J98").Value = 456.4
$ws.Range("K98").Value = 853
$ws.Range("L98").Value = 456.4
$ws.Range("M98").Value = 645
$ws.Range("N98").Value = -3452.4
$ws.Range("H106").Value = 7999.3335
$ws.Range("J106").Value = 7999.5
$ws.Range("L106").Value = 7999.5
$ws.Range("N106").Value = -9261.5
$ws.Range("H122").Value = 672.7273
$ws.Range("I122").Value = 853
$ws.Range("J122").Value = 456.4
$ws.Range("K122").Value = 2559
$ws.Range("L122").Value = 1369.2
$ws.Range("M122").Value = -109
$ws.Range("N122").Value = -6269.2
$ws.Range("H138").Value = 3148.0781
$ws.Range("J138").Value = 3229.0908
$ws.Range("L138").Value = 9687.2724
$ws.Range("N138").Value = -19967.2724

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1595.2858
$ws.Range("J4").Value = 1419.75
$ws.Range("L4").Value = 1419.75
$ws.Range("N4").Value = -1651.75
$ws.Range("H9").Value = 30009
$ws.Range("J9").Value = 30009
$ws.Range("L9").Value = 30009
$ws.Range("N9").Value = -30349
$ws.Range("H20").Value = 30009
$ws.Range("J20").Value = 30009
$ws.Range("L20").Value = 30009
$ws.Range("N20").Value = -30549
$ws.Range("H132").Value = 4332.5557
$ws.Range("I132").Value = 3798
$ws.Range("K132").Value = 11394
$ws.Range("M132").Value = -8864

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2503.3333
$ws.Range("I134").Value = 2107.6667
$ws.Range("K134").Value = 6323.000100000001
$ws.Range("M134").Value = -3788.000100000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2238.25
$ws.Range("I16").Value = 2238.25
$ws.Range("K16").Value = 2238.25
$ws.Range("M16").Value = -1951.25
$ws.Range("H59").Value = 60127
$ws.Range("J59").Value = 60127
$ws.Range("L59").Value = 60127
$ws.Range("N59").Value = -62417
$ws.Range("H68").Value = 55747.5
$ws.Range("I68").Value = 13000
$ws.Range("K68").Value = 13000
$ws.Range("M68").Value = -12251
$ws.Range("H71").Value = 55747.5
$ws.Range("I71").Value = 13000
$ws.Range("K71").Value = 39000
$ws.Range("M71").Value = -35256
$ws.Range("H74").Value = 69996.664
$ws.Range("J74").Value = 69996.664
$ws.Range("L74").Value = 69996.664
$ws.Range("N74").Value = -71744.664
$ws.Range("H77").Value = 69996.664
$ws.Range("J77").Value = 69996.664
$ws.Range("L77").Value = 209989.992
$ws.Range("N77").Value = -218725.992
$ws.Range("H107").Value = 2595.3845
$ws.Range("I107").Value = 999.5
$ws.Range("J107").Value = 2885.5454
$ws.Range("K107").Value = 999.5
$ws.Range("L107").Value = 2885.5454
$ws.Range("M107").Value = 920.5
$ws.Range("N107").Value = -6725.5454
$ws.Range("H113").Value = 2238.25
$ws.Range("I113").Value = 2238.25
$ws.Range("K113").Value = 2238.25
$ws.Range("M113").Value = -68.25
$ws.Range("H132").Value = 3068.4092
$ws.Range("I132").Value = 2599.25
$ws.Range("K132").Value = 7797.75
$ws.Range("M132").Value = -5267.75
$ws.Range("H134").Value = 3183.3333
$ws.Range("I134").Value = 1875
$ws.Range("J134").Value = 3837.5
$ws.Range("K134").Value = 5625
$ws.Range("L134").Value = 11512.5
$ws.Range("M134").Value = -3090
$ws.Range("N134").Value = -16582.5
$ws.Range("H141").Value = 61651.125
$ws.Range("J141").Value = 56172.715
$ws.Range("L141").Value = 56172.715
$ws.Range("N141").Value = -66532.715

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 1750
$ws.Range("I109").Value = 1750
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 5250
$ws.Range("L109").Value = 0
$ws.Range("M109").Value = -4210
$ws.Range("N109").ClearContents()
$ws.Range("H117").Value = 1516.3334
$ws.Range("I117").Value = 1125
$ws.Range("K117").Value = 3375
$ws.Range("M117").Value = 67
$ws.Range("H134").Value = 2880
$ws.Range("I134").Value = 2880
$ws.Range("K134").Value = 8640
$ws.Range("M134").Value = -3570
$ws.Range("H139").Value = 1026.8572
$ws.Range("I139").Value = 1031.3334
$ws.Range("J139").Value = 1000
$ws.Range("K139").Value = 3094.0002
$ws.Range("L139").Value = 3000
$ws.Range("M139").Value = 2045.9998
$ws.Range("N139").Value = -13280

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 88.46154
$ws.Range("I2").Value = 104
$ws.Range("J2").Value = 78.75
$ws.Range("K2").Value = 104
$ws.Range("L2").Value = 78.75
$ws.Range("M2").Value = 9
$ws.Range("N2").Value = -304.75
$ws.Range("H52").Value = 20000
$ws.Range("J52").Value = 20000
$ws.Range("L52").Value = 20000
$ws.Range("N52").Value = -20518
$ws.Range("H97").Value = 728.44446
$ws.Range("I97").Value = 987.8333
$ws.Range("K97").Value = 987.8333
$ws.Range("M97").Value = -491.8333
$ws.Range("H132").Value = 3889.7778
$ws.Range("I132").Value = 2337.3333
$ws.Range("K132").Value = 7011.999899999999
$ws.Range("M132").Value = -4481.999899999999
$ws.Range("H136").Value = 52875
$ws.Range("J136").Value = 52875
$ws.Range("L136").Value = 158625
$ws.Range("N136").Value = -163725

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 232.5
$ws.Range("I55").Value = 300.25
$ws.Range("J55").Value = 187.33333
$ws.Range("K55").Value = 300.25
$ws.Range("L55").Value = 187.33333
$ws.Range("M55").Value = -127.25
$ws.Range("N55").Value = -533.3333299999999

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 18503.5
$ws.Range("J15").Value = 18503.5
$ws.Range("L15").Value = 18503.5
$ws.Range("N15").Value = -19079.5
$ws.Range("H62").Value = 102499.5
$ws.Range("I62").Value = 102499.5
$ws.Range("K62").Value = 102499.5
$ws.Range("M62").Value = -101875.5
$ws.Range("H65").Value = 102499.5
$ws.Range("I65").Value = 102499.5
$ws.Range("K65").Value = 512497.5
$ws.Range("M65").Value = -509377.5
$ws.Range("H113").Value = 1857.1666
$ws.Range("I113").Value = 1830
$ws.Range("J113").Value = 1993
$ws.Range("K113").Value = 5490
$ws.Range("L113").Value = 5979
$ws.Range("M113").Value = -3320
$ws.Range("N113").Value = -10319
$ws.Range("H119").Value = 70000
$ws.Range("J119").Value = 70000
$ws.Range("L119").Value = 70000
$ws.Range("N119").Value = -79676
$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920
$ws.Range("H136").Value = 2464.261
$ws.Range("I136").Value = 1912.5714
$ws.Range("K136").Value = 5737.7142
$ws.Range("M136").Value = -3187.7142
